# "Drop in files from RMI script"
# Update the OCCF workbook's base-year dollar figures from 2018 to 2019:
#   - About!A26 literal conversion-factor value changes
#   - About! labels referencing "2018 dollars" (billion/million variants and
#     the "2012 dollars per 2018 dollar" note) become "2019 dollars"
#   - The dependent OCCF-Dp{L,M,S}OCU sheets recompute from the formulas that
#     already reference About!A26

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# New conversion-factor input value (About!A26). The three output sheets
# (OCCF-DpLOCU / OCCF-DpMOCU / OCCF-DpSOCU) hold formulas keyed off this cell
# and will recalculate automatically.
$ws.Range("A26").Value = 0.89805481563188172

# Relabel the unit captions / notes for the new 2019 base year. Order matches
# the order the strings were (re)written in the source workbook.
$ws.Range("B26").Value = "2019 dollars per 2012 dollar"
$ws.Range("B29").Value = 'which in this case is "2012 dollars per 2019 dollar."'
$ws.Range("A21").Value = "million 2019 dollars"
$ws.Range("A18").Value = "billion 2019 dollars"

# Restore the workbook's last-saved selection on the About sheet.
[void]$ws.Range("A19").Select()

# Turn on iterative calculation (workbook has circular-ish recalculation
# dependencies in the RMI model this feeds); mirrors the authored calcPr.
$excel.Iteration = $true
$excel.MaxIterations = 100
$excel.MaxChange = 0.00001
